$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.3300591674014615
$ws.Range("B2").Value = 0.3839125079094902
$ws.Range("C2").Value = 0.2823455332341271
$ws.Range("D2").Value = 0.3552559731737231
$ws.Range("E2").Value = 0.303890232326717
$ws.Range("K2").Value = 2.31041417181023
$ws.Range("L2").Value = 2.687387555366431
$ws.Range("M2").Value = 1.97641873263889
$ws.Range("N2").Value = 2.486791812216062
$ws.Range("O2").Value = 2.127231626287019
$ws.Range("P2").Value = 63.74049
$ws.Range("Q2").Value = 81.81617043041194
$ws.Range("R2").Value = 50.1719375637459
$ws.Range("S2").Value = 71.40129667373306
$ws.Range("T2").Value = 55.49574629892218
$ws.Range("U2").Value = 0.2037483695295493
$ws.Range("V2").Value = 0.2596457149369613
$ws.Range("W2").Value = 0.1490927975476136
$ws.Range("X2").Value = 0.2325214558897186
$ws.Range("Y2").Value = 0.1737636220520615
$ws.Range("Z2").Value = 0.8601408267920645
$ws.Range("AA2").Value = 0.915683397578995
$ws.Range("AB2").Value = 0.7913070969371712
$ws.Range("AC2").Value = 0.8905303427093396
$ws.Range("AD2").Value = 0.8250384108230511
$ws.Range("F3").Value = 7.000697589884063
$ws.Range("G3").Value = 8.841508045435681
$ws.Range("H3").Value = 5.460560086407628
$ws.Range("I3").Value = 7.836787123025934
$ws.Range("J3").Value = 6.132332668897923
$ws.Range("K3").Value = 2.310230204661741
$ws.Range("L3").Value = 2.917697654993776
$ws.Range("M3").Value = 1.801984828514517
$ws.Range("N3").Value = 2.586139750598559
$ws.Range("O3").Value = 2.023669780736315
$ws.Range("P3").Value = 63.7481
$ws.Range("Q3").Value = 75.27319761923111
$ws.Range("R3").Value = 56.19496914266413
$ws.Range("S3").Value = 68.56141804928605
$ws.Range("T3").Value = 59.10072488884368
$ws.Range("U3").Value = 0.2024537036706112
$ws.Range("V3").Value = 0.2897536462010484
$ws.Range("W3").Value = 0.1181684385509804
$ws.Range("X3").Value = 0.247716318828288
$ws.Range("Y3").Value = 0.1555729590431187
$ws.Range("Z3").Value = 0.8540245377345156
$ws.Range("AA3").Value = 0.9382565481519122
$ws.Range("AB3").Value = 0.7363604334518572
$ws.Range("AC3").Value = 0.9015636713450195
$ws.Range("AD3").Value = 0.7961231546180072
$ws.Range("A4").Value = 0.3299986630848838
$ws.Range("B4").Value = 0.383608226101057
$ws.Range("C4").Value = 0.2823486256541262
$ws.Range("D4").Value = 0.3549667219716015
$ws.Range("E4").Value = 0.3040319195507326
$ws.Range("F4").Value = 6.998388618331686
$ws.Range("G4").Value = 8.839401018561805
$ws.Range("H4").Value = 5.457049546651731
$ws.Range("I4").Value = 7.832553835399432
$ws.Range("J4").Value = 6.127876753459891
$ws.Range("K4").Value = 2.309357765733866
$ws.Range("L4").Value = 3.040405733971514
$ws.Range("M4").Value = 1.720088389927188
$ws.Range("N4").Value = 2.631383316563059
$ws.Range("O4").Value = 1.971117758847862
$ws.Range("P4").Value = 64.83936
$ws.Range("Q4").Value = 89.48542648455536
$ws.Range("R4").Value = 49.17162874104899
$ws.Range("S4").Value = 73.94453148713578
$ws.Range("T4").Value = 55.46041706019108
$ws.Range("U4").Value = 0.2015009989411887
$ws.Range("V4").Value = 0.3043881193105049
$ws.Range("W4").Value = 0.1028152436196753
$ws.Range("X4").Value = 0.254685267863865
$ws.Range("Y4").Value = 0.1465640520801581
$ws.Range("Z4").Value = 0.8497637062567935
$ws.Range("AA4").Value = 0.9483438917022717
$ws.Range("AB4").Value = 0.7036543722556876
$ws.Range("AC4").Value = 0.9057098087361612
$ws.Range("AD4").Value = 0.7796386398636745